$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Table cell margin: left padding 118 dxa (5.9 pt) -> 123 dxa (6.15 pt)
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.LeftPadding = 6.15

# ------------------------------------------------------------------
# 2. Merge the two runs holding the created-at date placeholder and the
#    trailing "г." into a single run whose text uses the "with quotes"
#    field variant.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "{issue.humanized_created_at_with_month_as_word} г.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{issue.humanized_created_at_with_quotes_and_month_as_word} г.",
    2) | Out-Null

# ------------------------------------------------------------------
# 3. Swap the names of the two co-located bookmark pairs (same start
#    position, so Delete + re-Add keeps the surrounding w:id sequence
#    intact while exchanging the w:name values).
# ------------------------------------------------------------------
$bmA = $d.Bookmarks.Item("_Hlk505353486")
$rangeA = $bmA.Range
$bmB = $d.Bookmarks.Item("_Hlk505353748")
$rangeB = $bmB.Range
$bmA.Delete()
$bmB.Delete()
$d.Bookmarks.Add("_Hlk505353748", $rangeA) | Out-Null
$d.Bookmarks.Add("_Hlk505353486", $rangeB) | Out-Null

$bmC = $d.Bookmarks.Item("_Hlk505353774")
$rangeC = $bmC.Range
$bmD = $d.Bookmarks.Item("_Hlk505353480")
$rangeD = $bmD.Range
$bmC.Delete()
$bmD.Delete()
$d.Bookmarks.Add("_Hlk505353480", $rangeC) | Out-Null
$d.Bookmarks.Add("_Hlk505353774", $rangeD) | Out-Null

# ------------------------------------------------------------------
# 4. Add the two new character styles used for the formatted list
#    labels ("ListLabel 9" bold / "ListLabel 10" non-bold).
# ------------------------------------------------------------------
$ll9 = $d.Styles.Add("ListLabel9", 2)
$ll9.NameLocal = "ListLabel 9"
$ll9.QuickStyle = $true
$ll9.Font.Name = "Times New Roman"
$ll9.Font.Bold = $true
$ll9.Font.Size = 10.5

$ll10 = $d.Styles.Add("ListLabel10", 2)
$ll10.NameLocal = "ListLabel 10"
$ll10.QuickStyle = $true
$ll10.Font.Name = "Times New Roman"
$ll10.Font.Bold = $false
$ll10.Font.Size = 10.5

Write-Output "done"
